$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.423.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.939.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.00%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "477.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +8.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.71"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.02%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.732"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000349"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +10.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.42"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.574.18"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.93"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.950.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.65%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.609.55"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.74"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.38"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.52"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.52"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.68"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.61"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "718.78"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.133"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.44"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.23%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.14"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0836"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +23.62%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.25"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.33%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.07%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.38"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.06"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.12%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.28%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.337"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.87%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.49"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.81%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.53"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.28"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.78"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.88"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.49"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.08%  "
